$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New calendar values for rows 2-6 (columns A-G)
$values = @(
    @("27/11","28/11","29/11","30/11","1/12","2/12","3/12"),
    @("4/12","5/12","6/12","7/12","8/12","9/12","10/12"),
    @("11/12","12/12","13/12","14/12","15/12","16/12","17/12"),
    @("18/12","19/12","20/12","21/12","22/12","23/12","24/12"),
    @("25/12","26/12","27/12","28/12","29/12","30/12","31/12")
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

# Remove old row 7 contents entirely, since the calendar now only spans through row 6
$ws.Range("A7:G7").Clear()
